# Generate Report for Handback
# Update handback-status report with refreshed generation timestamps and
# priority values for the 2b0f4cc5... / 94790c54... rows.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" for the 2b0f4cc5... / 94790c54... rows
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-01 16:19:59"
$wsOverview.Range("G4").Value = "2016-09-01 16:19:59"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-09-01 16:19:54"
$wsZhCn.Range("H4").Value = "2016-09-01 16:19:54"
$wsZhCn.Range("K3").Value = "2016-09-01 16:20:40"
$wsZhCn.Range("K4").Value = "2016-09-01 16:20:40"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-09-01 16:19:59"
$wsDeDe.Range("H4").Value = "2016-09-01 16:19:59"
$wsDeDe.Range("K3").Value = "2016-09-01 16:20:47"
$wsDeDe.Range("K4").Value = "2016-09-01 16:20:47"
